# Auto-generated Excel COM-interop script to apply the Cactuar_Profits value updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4166.6665
$ws.Range("J64").Value = 4000
$ws.Range("L64").Value = 4000
$ws.Range("N64").Value = -4496
$ws.Range("H67").Value = 4166.6665
$ws.Range("J67").Value = 4000
$ws.Range("L67").Value = 4000
$ws.Range("N67").Value = -5716
$ws.Range("H137").Value = 9285289
$ws.Range("I137").Value = 460034.2
$ws.Range("J137").Value = 22228996
$ws.Range("K137").Value = 1380102.6
$ws.Range("L137").Value = 66686988
$ws.Range("M137").Value = -1377552.6
$ws.Range("N137").Value = -66692088
$ws.Range("H138").Value = 7475.154
$ws.Range("I138").Value = 5095.5
$ws.Range("J138").Value = 7907.8184
$ws.Range("K138").Value = 15286.5
$ws.Range("L138").Value = 23723.4552
$ws.Range("M138").Value = -10146.5
$ws.Range("N138").Value = -34003.4552

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 17597.334
$ws.Range("I2").Value = 20023.846
$ws.Range("K2").Value = 20023.846
$ws.Range("M2").Value = -19910.846
$ws.Range("H45").Value = 41357.484
$ws.Range("I45").Value = 58545.45
$ws.Range("K45").Value = 58545.45
$ws.Range("M45").Value = -58168.45
$ws.Range("H61").Value = 3814.59
$ws.Range("I61").Value = 3493.7856
$ws.Range("J61").Value = 4523.737
$ws.Range("K61").Value = 3493.7856
$ws.Range("L61").Value = 4523.737
$ws.Range("M61").Value = -3281.7856
$ws.Range("N61").Value = -4947.737
$ws.Range("H63").Value = 2069.5715
$ws.Range("I63").Value = 1729
$ws.Range("K63").Value = 1729
$ws.Range("M63").Value = -1043
$ws.Range("H66").Value = 2069.5715
$ws.Range("I66").Value = 1729
$ws.Range("K66").Value = 8645
$ws.Range("M66").Value = -5213
$ws.Range("H74").Value = 1387.2858
$ws.Range("I74").Value = 1165.2858
$ws.Range("K74").Value = 1165.2858
$ws.Range("M74").Value = -291.2858000000001
$ws.Range("H77").Value = 1387.2858
$ws.Range("I77").Value = 1165.2858
$ws.Range("K77").Value = 5826.429
$ws.Range("M77").Value = -1458.429
$ws.Range("H102").Value = 1703.5834
$ws.Range("I102").Value = 1683
$ws.Range("J102").Value = 1765.3334
$ws.Range("K102").Value = 1683
$ws.Range("L102").Value = 1765.3334
$ws.Range("M102").Value = -61
$ws.Range("N102").Value = -5009.3334
$ws.Range("H110").Value = 1217.9131
$ws.Range("I110").Value = 980.6667
$ws.Range("K110").Value = 980.6667
$ws.Range("M110").Value = 1064.3333
$ws.Range("H116").Value = 17597.334
$ws.Range("I116").Value = 20023.846
$ws.Range("K116").Value = 20023.846
$ws.Range("M116").Value = -17729.846
$ws.Range("H122").Value = 4547.421
$ws.Range("I122").Value = 2918.2727
$ws.Range("K122").Value = 8754.8181
$ws.Range("M122").Value = -6304.8181
$ws.Range("H136").Value = 3814.59
$ws.Range("I136").Value = 3493.7856
$ws.Range("J136").Value = 4523.737
$ws.Range("K136").Value = 10481.3568
$ws.Range("L136").Value = 13571.211
$ws.Range("M136").Value = -7931.356800000001
$ws.Range("N136").Value = -18671.211

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 17597.334
$ws.Range("I3").Value = 20023.846
$ws.Range("K3").Value = 20023.846
$ws.Range("M3").Value = -19909.846
$ws.Range("H26").Value = 10246
$ws.Range("I26").Value = 10246
$ws.Range("K26").Value = 10246
$ws.Range("M26").Value = -9954
$ws.Range("H86").Value = 4101.9546
$ws.Range("J86").Value = 4543.5
$ws.Range("L86").Value = 4543.5
$ws.Range("N86").Value = -6789.5
$ws.Range("H89").Value = 4101.9546
$ws.Range("J89").Value = 4543.5
$ws.Range("L89").Value = 22717.5
$ws.Range("N89").Value = -33949.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 300500
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 300500
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 300500
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -300780
$ws.Range("H31").Value = 25643210
$ws.Range("I31").Value = 28573048
$ws.Range("K31").Value = 28573048
$ws.Range("M31").Value = -28572753
$ws.Range("H34").Value = 25643210
$ws.Range("I34").Value = 28573048
$ws.Range("K34").Value = 28573048
$ws.Range("M34").Value = -28572846
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H132").Value = 60608310
$ws.Range("I132").Value = 70177050
$ws.Range("K132").Value = 210531150
$ws.Range("M132").Value = -210528620
$ws.Range("H141").Value = 201543.14
$ws.Range("J141").Value = 217319.08
$ws.Range("L141").Value = 217319.08
$ws.Range("N141").Value = -227679.08

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 257.58334
$ws.Range("I14").Value = 257.58334
$ws.Range("K14").Value = 772.7500200000001
$ws.Range("M14").Value = -599.7500200000001
$ws.Range("H34").Value = 1243.4445
$ws.Range("I34").Value = 353.36365
$ws.Range("K34").Value = 1060.09095
$ws.Range("M34").Value = -976.09095
$ws.Range("H80").Value = 5988.6
$ws.Range("J80").Value = 5988.6
$ws.Range("L80").Value = 17965.8
$ws.Range("N80").Value = -19837.8
$ws.Range("H83").Value = 5988.6
$ws.Range("J83").Value = 5988.6
$ws.Range("L83").Value = 53897.4
$ws.Range("N83").Value = -63257.4
$ws.Range("H122").Value = 1506.125
$ws.Range("I122").Value = 1106.4
$ws.Range("K122").Value = 9957.6
$ws.Range("M122").Value = -7507.6
$ws.Range("H131").Value = 17117916
$ws.Range("J131").Value = 19938704
$ws.Range("L131").Value = 59816112
$ws.Range("N131").Value = -59826192

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5680.524
$ws.Range("J70").Value = 5942.1816
$ws.Range("L70").Value = 5942.1816
$ws.Range("N70").Value = -6482.1816
$ws.Range("H73").Value = 5680.524
$ws.Range("J73").Value = 5942.1816
$ws.Range("L73").Value = 5942.1816
$ws.Range("N73").Value = -7814.1816
$ws.Range("H123").Value = 43336.4
$ws.Range("J123").Value = 43336.4
$ws.Range("L123").Value = 43336.4
$ws.Range("N123").Value = -48236.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1061.4615
$ws.Range("J22").Value = 1487.5
$ws.Range("L22").Value = 1487.5
$ws.Range("N22").Value = -2077.5
$ws.Range("H27").Value = 1061.4615
$ws.Range("J27").Value = 1487.5
$ws.Range("L27").Value = 1487.5
$ws.Range("N27").Value = -1701.5
$ws.Range("H30").Value = 745
$ws.Range("I30").Value = 772
$ws.Range("J30").Value = 718
$ws.Range("K30").Value = 772
$ws.Range("L30").Value = 718
$ws.Range("M30").Value = -664
$ws.Range("N30").Value = -934
$ws.Range("H31").Value = 2524.5
$ws.Range("I31").Value = 448
$ws.Range("K31").Value = 448
$ws.Range("M31").Value = -200
$ws.Range("H136").Value = 5211.9
$ws.Range("I136").Value = 3024.739
$ws.Range("K136").Value = 9074.217000000001
$ws.Range("M136").Value = -6524.217000000001
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 8062.8237
$ws.Range("I81").Value = 6459.25
$ws.Range("K81").Value = 12918.5
$ws.Range("M81").Value = -11857.5
$ws.Range("H84").Value = 8062.8237
$ws.Range("I84").Value = 6459.25
$ws.Range("K84").Value = 64592.5
$ws.Range("M84").Value = -59288.5
$ws.Range("H93").Value = 39959.25
$ws.Range("J93").Value = 39959.25
$ws.Range("L93").Value = 39959.25
$ws.Range("N93").Value = -44951.25
$ws.Range("H107").Value = 949.0833
$ws.Range("I107").Value = 782
$ws.Range("J107").Value = 1116.1666
$ws.Range("K107").Value = 2346
$ws.Range("L107").Value = 3348.4998
$ws.Range("M107").Value = -426
$ws.Range("N107").Value = -7188.4998
$ws.Range("H132").Value = 3173.75
$ws.Range("I132").Value = 1515.75
$ws.Range("K132").Value = 4547.25
$ws.Range("M132").Value = -2017.25
$ws.Range("H136").Value = 5051.5625
$ws.Range("I136").Value = 3562.238
$ws.Range("K136").Value = 10686.714
$ws.Range("M136").Value = -8136.714
$ws.Range("H139").Value = 175000
$ws.Range("J139").Value = 230000
$ws.Range("L139").Value = 230000
$ws.Range("N139").Value = -240280
$ws.Range("H140").Value = 99499.5
$ws.Range("J140").Value = 99499.5
$ws.Range("L140").Value = 99499.5
$ws.Range("N140").Value = -109859.5
$ws.Range("H141").Value = 61875
$ws.Range("J141").Value = 61875
$ws.Range("L141").Value = 61875
$ws.Range("N141").Value = -72235
